$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.748.65"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.629.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.65%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.500"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.57%  "
$ws.Range("E8").Value = "  -0.87%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0632"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.51"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0791"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.02%  "
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.855.99"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.628.17"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.552"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₃0759"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.85"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.748.72"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("E21").Value = "  -0.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.90"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.79%  "
$ws.Range("B24").Value = "BinanceUSD"
$ws.Range("C24").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.61%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.82"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.64%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.50"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.62%  "
$ws.Range("E27").Value = "  +3.83%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.83"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.32%  "
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("E30").Value = "  -0.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0494"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.23"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.54%  "
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("E35").Value = "  -0.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.907"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.139.54"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.96%  "
$ws.Range("E38").Value = "  -2.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.542"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0155"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.38%  "
$ws.Range("E41").Value = "  -0.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.52"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.55"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.803"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.763.85"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.23%  "
$ws.Range("E48").Value = "  +1.66%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.417"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.42%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.44"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.33"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.71%  "
